$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before the current row 69. This shifts the
# existing rows 69-82 down to 71-84 (dimension grows from A1:T82 to A1:T84).
$ws.Rows("69:70").Insert()

# ---- New row 69 -----------------------------------------------------
$ws.Cells.Item(69, 1).Value = 10
$ws.Cells.Item(69, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(69, 3).Value = "La Araucanía"
$ws.Cells.Item(69, 4).Value = 44559
$ws.Cells.Item(69, 5).Value = 9
$ws.Cells.Item(69, 6).Value = "Fruta"
$ws.Cells.Item(69, 7).Value = 100101
$ws.Cells.Item(69, 8).Value = "Berries"
$ws.Cells.Item(69, 9).Value = 100101001
$ws.Cells.Item(69, 10).Value = "Arándano (blue)"
$ws.Cells.Item(69, 11).Value = "Sin especificar"
$ws.Cells.Item(69, 12).Value = "Primera"
$ws.Cells.Item(69, 13).Value = 300
$ws.Cells.Item(69, 14).Value = 1500
$ws.Cells.Item(69, 15).Value = 1500
$ws.Cells.Item(69, 16).Value = 1500
$ws.Cells.Item(69, 17).Value = "$/kilo"
$ws.Cells.Item(69, 18).Value = "Región de La Araucanía"
$ws.Cells.Item(69, 19).Value = 1500
$ws.Cells.Item(69, 20).Value = 1

# ---- New row 70 -----------------------------------------------------
$ws.Cells.Item(70, 1).Value = 10
$ws.Cells.Item(70, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(70, 3).Value = "La Araucanía"
$ws.Cells.Item(70, 4).Value = 44559
$ws.Cells.Item(70, 5).Value = 9
$ws.Cells.Item(70, 6).Value = "Fruta"
$ws.Cells.Item(70, 7).Value = 100101
$ws.Cells.Item(70, 8).Value = "Berries"
$ws.Cells.Item(70, 9).Value = 100101001
$ws.Cells.Item(70, 10).Value = "Arándano (blue)"
$ws.Cells.Item(70, 11).Value = "Sin especificar"
$ws.Cells.Item(70, 12).Value = "Primera"
$ws.Cells.Item(70, 13).Value = 200
$ws.Cells.Item(70, 14).Value = 2200
$ws.Cells.Item(70, 15).Value = 2200
$ws.Cells.Item(70, 16).Value = 2200
$ws.Cells.Item(70, 17).Value = "$/kilo"
$ws.Cells.Item(70, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(70, 19).Value = 2200
$ws.Cells.Item(70, 20).Value = 1
